$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.590955666666667
$ws.Range("H2").Value = 4.772867
$ws.Range("I2").Value = 0.4719498906366805
$ws.Range("J2").Value = 0.4719498906366804
$ws.Range("M2").Value = 68.18146900000001
$ws.Range("N2").Value = 204.544407
$ws.Range("O2").Value = 0.1244286043321187
$ws.Range("P2").Value = 0.1244286043321187
$ws.Range("Q2").Value = 108.4736944672077
$ws.Range("R2").Value = 976.263250204869
$ws.Range("S2").Value = 0.05872406620661819
$ws.Range("T2").Value = 0.05872406620661818
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.590955666666667
$ws.Range("H3").Value = 4.772867
$ws.Range("I3").Value = 0.4719498906366805
$ws.Range("J3").Value = 0.4719498906366804
$ws.Range("O3").Value = 0.345973452289334
$ws.Range("P3").Value = 0.3459734522893341
$ws.Range("Q3").Value = 301.6108615767133
$ws.Range("R3").Value = 2714.49775419042
$ws.Range("S3").Value = 0.163282132971146
$ws.Range("T3").Value = 0.163282132971146
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.590955666666667
$ws.Range("H4").Value = 4.772867
$ws.Range("I4").Value = 0.4719498906366805
$ws.Range("J4").Value = 0.4719498906366804
$ws.Range("M4").Value = 188.0130056666667
$ws.Range("N4").Value = 564.0390170000001
$ws.Range("O4").Value = 0.3431166302883566
$ws.Range("P4").Value = 0.3431166302883567
$ws.Range("Q4").Value = 299.1203567724154
$ws.Range("R4").Value = 2692.083210951739
$ws.Range("S4").Value = 0.1619338561402162
$ws.Range("T4").Value = 0.1619338561402162
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.590955666666667
$ws.Range("H5").Value = 4.772867
$ws.Range("I5").Value = 0.4719498906366805
$ws.Range("J5").Value = 0.4719498906366804
$ws.Range("M5").Value = 102.1836573333333
$ws.Range("N5").Value = 306.550972
$ws.Range("O5").Value = 0.1864813130901906
$ws.Range("P5").Value = 0.1864813130901907
$ws.Range("Q5").Value = 162.5696686751915
$ws.Range("R5").Value = 1463.127018076724
$ws.Range("S5").Value = 0.08800983531870005
$ws.Range("T5").Value = 0.08800983531870005
$ws.Range("I6").Value = 0.5280501093633195
$ws.Range("J6").Value = 0.5280501093633195
$ws.Range("M6").Value = 68.18146900000001
$ws.Range("N6").Value = 204.544407
$ws.Range("O6").Value = 0.1244286043321187
$ws.Range("P6").Value = 0.1244286043321187
$ws.Range("Q6").Value = 121.367855704299
$ws.Range("R6").Value = 1092.310701338691
$ws.Range("S6").Value = 0.06570453812550048
$ws.Range("T6").Value = 0.06570453812550048
$ws.Range("I7").Value = 0.5280501093633195
$ws.Range("J7").Value = 0.5280501093633195
$ws.Range("O7").Value = 0.345973452289334
$ws.Range("P7").Value = 0.3459734522893341
$ws.Range("S7").Value = 0.1826913193181881
$ws.Range("T7").Value = 0.1826913193181881
$ws.Range("I8").Value = 0.5280501093633195
$ws.Range("J8").Value = 0.5280501093633195
$ws.Range("M8").Value = 188.0130056666667
$ws.Range("N8").Value = 564.0390170000001
$ws.Range("O8").Value = 0.3431166302883566
$ws.Range("P8").Value = 0.3431166302883567
$ws.Range("Q8").Value = 334.676499010069
$ws.Range("R8").Value = 3012.088491090621
$ws.Range("S8").Value = 0.1811827741481404
$ws.Range("T8").Value = 0.1811827741481404
$ws.Range("I9").Value = 0.5280501093633195
$ws.Range("J9").Value = 0.5280501093633195
$ws.Range("M9").Value = 102.1836573333333
$ws.Range("N9").Value = 306.550972
$ws.Range("O9").Value = 0.1864813130901906
$ws.Range("P9").Value = 0.1864813130901907
$ws.Range("Q9").Value = 181.894165093004
$ws.Range("R9").Value = 1637.047485837036
$ws.Range("S9").Value = 0.0984714777714906
$ws.Range("T9").Value = 0.09847147777149061
